$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# The B2:B3 merge is being split into two distinct "tipoDemonstracao" rows,
# so the merge must go away before the cells can hold independent values.
$ws1.Range("B2:B3").UnMerge()

# B2 and B3 move off the centered-alignment style onto the new "general
# horizontal / top vertical" style (new cellXfs entry).
$ws1.Range("B2").HorizontalAlignment = 1
$ws1.Range("B3").HorizontalAlignment = 1
$ws1.Range("B3").VerticalAlignment = -4160

# New statement-type values: B3 becomes "Balanço Patrimonial (BP)" and B4
# becomes "Demonstração do Resultado do Exercício (DRE)" (B2 keeps DCC).
$ws1.Range("B3").Value2 = "Balanço Patrimonial (BP)"
$ws1.Range("B4").Value2 = "Demonstração do Resultado do Exercício (DRE)"

# "demonstracoes" becomes the active/selected sheet with B5 selected;
# "cnpjs" is no longer the tab shown when the workbook opens.
$ws1.Activate()
$ws1.Range("B5").Select()
